$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A8:AB8").EntireRow.RowHeight = 54
